$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1283
$ws.Cells.Item(1283, 1).Value = "cs"
$ws.Cells.Item(1283, 1).WrapText = $true
$ws.Cells.Item(1283, 1).Font.Size = 10
$ws.Cells.Item(1283, 2).Value = "market.license.label"
$ws.Cells.Item(1283, 2).WrapText = $true
$ws.Cells.Item(1283, 2).Font.Size = 10
$ws.Cells.Item(1283, 3).Value = "Licence"
$ws.Cells.Item(1283, 3).WrapText = $true
$ws.Cells.Item(1283, 3).Font.Size = 10

# Row 1284
$ws.Cells.Item(1284, 1).Value = "cs"
$ws.Cells.Item(1284, 1).WrapText = $true
$ws.Cells.Item(1284, 1).Font.Size = 10
$ws.Cells.Item(1284, 2).Value = "market.license.index.title"
$ws.Cells.Item(1284, 2).WrapText = $true
$ws.Cells.Item(1284, 2).Font.Size = 10
$ws.Cells.Item(1284, 3).Value = "Licence"
$ws.Cells.Item(1284, 3).WrapText = $true
$ws.Cells.Item(1284, 3).Font.Size = 10

# Row 1285
$ws.Cells.Item(1285, 1).Value = "cs"
$ws.Cells.Item(1285, 1).WrapText = $true
$ws.Cells.Item(1285, 1).Font.Size = 10
$ws.Cells.Item(1285, 2).Value = "market.license.create.button"
$ws.Cells.Item(1285, 2).WrapText = $true
$ws.Cells.Item(1285, 2).Font.Size = 10
$ws.Cells.Item(1285, 3).Value = "Nová licence"
$ws.Cells.Item(1285, 3).WrapText = $true
$ws.Cells.Item(1285, 3).Font.Size = 10

# Row 1286
$ws.Cells.Item(1286, 1).Value = "cs"
$ws.Cells.Item(1286, 1).WrapText = $true
$ws.Cells.Item(1286, 1).Font.Size = 10
$ws.Cells.Item(1286, 2).Value = "market.license.create.title"
$ws.Cells.Item(1286, 2).WrapText = $true
$ws.Cells.Item(1286, 2).Font.Size = 10
$ws.Cells.Item(1286, 3).Value = "Nová licence"
$ws.Cells.Item(1286, 3).WrapText = $true
$ws.Cells.Item(1286, 3).Font.Size = 10

# Row 1287
$ws.Cells.Item(1287, 1).Value = "cs"
$ws.Cells.Item(1287, 1).WrapText = $true
$ws.Cells.Item(1287, 1).Font.Size = 10
$ws.Cells.Item(1287, 2).Value = "shared.license.create.name.label"
$ws.Cells.Item(1287, 2).WrapText = $true
$ws.Cells.Item(1287, 2).Font.Size = 10
$ws.Cells.Item(1287, 3).Value = "Název licence"
$ws.Cells.Item(1287, 3).WrapText = $true
$ws.Cells.Item(1287, 3).Font.Size = 10

# Row 1288
$ws.Cells.Item(1288, 1).Value = "cs"
$ws.Cells.Item(1288, 1).WrapText = $true
$ws.Cells.Item(1288, 1).Font.Size = 10
$ws.Cells.Item(1288, 2).Value = "shared.license.create.code.label"
$ws.Cells.Item(1288, 2).WrapText = $true
$ws.Cells.Item(1288, 2).Font.Size = 10
$ws.Cells.Item(1288, 3).Value = "Kód licence"
$ws.Cells.Item(1288, 3).WrapText = $true
$ws.Cells.Item(1288, 3).Font.Size = 10

# Row 1289
$ws.Cells.Item(1289, 1).Value = "cs"
$ws.Cells.Item(1289, 1).WrapText = $true
$ws.Cells.Item(1289, 1).Font.Size = 10
$ws.Cells.Item(1289, 2).Value = "shared.license.create.cost.label"
$ws.Cells.Item(1289, 2).WrapText = $true
$ws.Cells.Item(1289, 2).Font.Size = 10
$ws.Cells.Item(1289, 3).Value = "Cena licence"
$ws.Cells.Item(1289, 3).WrapText = $true
$ws.Cells.Item(1289, 3).Font.Size = 10

# Row 1290
$ws.Cells.Item(1290, 1).Value = "cs"
$ws.Cells.Item(1290, 1).WrapText = $true
$ws.Cells.Item(1290, 1).Font.Size = 10
$ws.Cells.Item(1290, 2).Value = "shared.license.create.renew.label"
$ws.Cells.Item(1290, 2).WrapText = $true
$ws.Cells.Item(1290, 2).Font.Size = 10
$ws.Cells.Item(1290, 3).Value = "Cena obnovy licence"
$ws.Cells.Item(1290, 3).WrapText = $true
$ws.Cells.Item(1290, 3).Font.Size = 10

# Row 1291
$ws.Cells.Item(1291, 1).Value = "cs"
$ws.Cells.Item(1291, 1).WrapText = $true
$ws.Cells.Item(1291, 1).Font.Size = 10
$ws.Cells.Item(1291, 2).Value = "shared.license.create.duration.label"
$ws.Cells.Item(1291, 2).WrapText = $true
$ws.Cells.Item(1291, 2).Font.Size = 10
$ws.Cells.Item(1291, 3).Value = "Platnost licence"
$ws.Cells.Item(1291, 3).WrapText = $true
$ws.Cells.Item(1291, 3).Font.Size = 10

# Row 1292
$ws.Cells.Item(1292, 1).Value = "cs"
$ws.Cells.Item(1292, 1).WrapText = $true
$ws.Cells.Item(1292, 1).Font.Size = 10
$ws.Cells.Item(1292, 2).Value = "shared.license.create.tokens.label"
$ws.Cells.Item(1292, 2).WrapText = $true
$ws.Cells.Item(1292, 2).Font.Size = 10
$ws.Cells.Item(1292, 3).Value = "Tokeny licence"
$ws.Cells.Item(1292, 3).WrapText = $true
$ws.Cells.Item(1292, 3).Font.Size = 10

# Row 1293
$ws.Cells.Item(1293, 1).Value = "cs"
$ws.Cells.Item(1293, 1).WrapText = $true
$ws.Cells.Item(1293, 1).Font.Size = 10
$ws.Cells.Item(1293, 2).Value = "shared.license.create.create"
$ws.Cells.Item(1293, 2).WrapText = $true
$ws.Cells.Item(1293, 2).Font.Size = 10
$ws.Cells.Item(1293, 3).Value = "Vytvořit licenci"
$ws.Cells.Item(1293, 3).WrapText = $true
$ws.Cells.Item(1293, 3).Font.Size = 10

# Row 1294
$ws.Cells.Item(1294, 1).Value = "cs"
$ws.Cells.Item(1294, 1).WrapText = $true
$ws.Cells.Item(1294, 1).Font.Size = 10
$ws.Cells.Item(1294, 2).Value = "shared.license.create.name.label.tooltip"
$ws.Cells.Item(1294, 2).WrapText = $true
$ws.Cells.Item(1294, 2).Font.Size = 10
$ws.Cells.Item(1294, 3).Value = "Použijte rozumné jméno pro licenci aby dávalo smysl; procházi překladem, tudíž je možné ho případně upravit jazykovou mutací."
$ws.Cells.Item(1294, 3).WrapText = $true
$ws.Cells.Item(1294, 3).Font.Size = 10

# Row 1295
$ws.Cells.Item(1295, 1).Value = "cs"
$ws.Cells.Item(1295, 1).WrapText = $true
$ws.Cells.Item(1295, 1).Font.Size = 10
$ws.Cells.Item(1295, 2).Value = "shared.license.create.code.label.tooltip"
$ws.Cells.Item(1295, 2).WrapText = $true
$ws.Cells.Item(1295, 2).Font.Size = 10
$ws.Cells.Item(1295, 3).Value = "Kód je unikátní označení, hlavně určený pro import/export."
$ws.Cells.Item(1295, 3).WrapText = $true
$ws.Cells.Item(1295, 3).Font.Size = 10

# Row 1296
$ws.Cells.Item(1296, 1).Value = "cs"
$ws.Cells.Item(1296, 1).WrapText = $true
$ws.Cells.Item(1296, 1).Font.Size = 10
$ws.Cells.Item(1296, 2).Value = "shared.license.create.cost.label.tooltip"
$ws.Cells.Item(1296, 2).WrapText = $true
$ws.Cells.Item(1296, 2).Font.Size = 10
$ws.Cells.Item(1296, 3).Value = "Pokud je uvedena cena licence, bude dostupná na tržišti."
$ws.Cells.Item(1296, 3).WrapText = $true
$ws.Cells.Item(1296, 3).Font.Size = 10

# Row 1297
$ws.Cells.Item(1297, 1).Value = "cs"
$ws.Cells.Item(1297, 1).WrapText = $true
$ws.Cells.Item(1297, 1).Font.Size = 10
$ws.Cells.Item(1297, 2).Value = "shared.license.create.renew.label.tooltip"
$ws.Cells.Item(1297, 2).WrapText = $true
$ws.Cells.Item(1297, 2).Font.Size = 10
$ws.Cells.Item(1297, 3).Value = "Pokud někdo již licenci vlastní a vyprší mu, může zaplatit uvedenou částku za obnovu; pokud není uvedena, obnova nebude možná."
$ws.Cells.Item(1297, 3).WrapText = $true
$ws.Cells.Item(1297, 3).Font.Size = 10

# Row 1298
$ws.Cells.Item(1298, 1).Value = "cs"
$ws.Cells.Item(1298, 1).WrapText = $true
$ws.Cells.Item(1298, 1).Font.Size = 10
$ws.Cells.Item(1298, 2).Value = "shared.license.create.duration.label.tooltip"
$ws.Cells.Item(1298, 2).WrapText = $true
$ws.Cells.Item(1298, 2).Font.Size = 10
$ws.Cells.Item(1298, 3).Value = "Délka trvání přidelené licence ve dnech od data jejího pořízení; licence obecně mají krátkodobý efekt."
$ws.Cells.Item(1298, 3).WrapText = $true
$ws.Cells.Item(1298, 3).Font.Size = 10

# Row 1299
$ws.Cells.Item(1299, 1).Value = "cs"
$ws.Cells.Item(1299, 1).WrapText = $true
$ws.Cells.Item(1299, 1).Font.Size = 10
$ws.Cells.Item(1299, 2).Value = "shared.license.create.tokens.label.tooltip"
$ws.Cells.Item(1299, 2).WrapText = $true
$ws.Cells.Item(1299, 2).Font.Size = 10
$ws.Cells.Item(1299, 3).Value = "Tokeny přidelené touto licencí. Prakticky udává oprávnění, jaká uživatel s touto licencí může získat."
$ws.Cells.Item(1299, 3).WrapText = $true
$ws.Cells.Item(1299, 3).Font.Size = 10

# Row 1300
$ws.Cells.Item(1300, 1).Value = "cs"
$ws.Cells.Item(1300, 1).WrapText = $true
$ws.Cells.Item(1300, 1).Font.Size = 10
$ws.Cells.Item(1300, 2).Value = "shared.license.create.success"
$ws.Cells.Item(1300, 2).WrapText = $true
$ws.Cells.Item(1300, 2).Font.Size = 10
$ws.Cells.Item(1300, 3).Value = "Licence [{{name}}] byla úspěšně vytvořena."
$ws.Cells.Item(1300, 3).WrapText = $true
$ws.Cells.Item(1300, 3).Font.Size = 10

# Row 1301
$ws.Cells.Item(1301, 1).Value = "cs"
$ws.Cells.Item(1301, 1).WrapText = $true
$ws.Cells.Item(1301, 1).Font.Size = 10
$ws.Cells.Item(1301, 2).Value = "market.license.delete.modal.title"
$ws.Cells.Item(1301, 2).WrapText = $true
$ws.Cells.Item(1301, 2).Font.Size = 10
$ws.Cells.Item(1301, 3).Value = "Odstranit vybrané licence?"
$ws.Cells.Item(1301, 3).WrapText = $true
$ws.Cells.Item(1301, 3).Font.Size = 10

# Row 1302
$ws.Cells.Item(1302, 1).Value = "cs"
$ws.Cells.Item(1302, 1).WrapText = $true
$ws.Cells.Item(1302, 1).Font.Size = 10
$ws.Cells.Item(1302, 2).Value = "market.license.delete.modal.content"
$ws.Cells.Item(1302, 2).WrapText = $true
$ws.Cells.Item(1302, 2).Font.Size = 10
$ws.Cells.Item(1302, 3).Value = "Tato akce obecně není tak destruktivní, jako smazání certifikátu, přesto se ale třikrát a půlkrát rozmyslete, co děláte, protože můžete naštvat hodně lidí."
$ws.Cells.Item(1302, 3).WrapText = $true
$ws.Cells.Item(1302, 3).Font.Size = 10

# Row 1303
$ws.Cells.Item(1303, 1).Value = "cs"
$ws.Cells.Item(1303, 1).WrapText = $true
$ws.Cells.Item(1303, 1).Font.Size = 10
$ws.Cells.Item(1303, 2).Value = "market.license.delete.success"
$ws.Cells.Item(1303, 2).WrapText = $true
$ws.Cells.Item(1303, 2).Font.Size = 10
$ws.Cells.Item(1303, 3).Value = "Vybrané licence byly úspěšně odstraněny."
$ws.Cells.Item(1303, 3).WrapText = $true
$ws.Cells.Item(1303, 3).Font.Size = 10

# Row 1304
$ws.Cells.Item(1304, 1).Value = "cs"
$ws.Cells.Item(1304, 1).WrapText = $true
$ws.Cells.Item(1304, 1).Font.Size = 10
$ws.Cells.Item(1304, 2).Value = "market.license.index.subtitle"
$ws.Cells.Item(1304, 2).WrapText = $true
$ws.Cells.Item(1304, 2).Font.Size = 10
$ws.Cells.Item(1304, 3).Value = "Krátkodobé propůjčení privilegií."
$ws.Cells.Item(1304, 3).WrapText = $true
$ws.Cells.Item(1304, 3).Font.Size = 10

# Row 1305
$ws.Cells.Item(1305, 1).Value = "cs"
$ws.Cells.Item(1305, 1).WrapText = $true
$ws.Cells.Item(1305, 1).Font.Size = 10
$ws.Cells.Item(1305, 2).Value = "market.license.index.content"
$ws.Cells.Item(1305, 2).WrapText = $true
$ws.Cells.Item(1305, 2).Font.Size = 10
$ws.Cells.Item(1305, 3).Value = "<p>`n`tLicence jsou slabší a chudější verzí Certifikátů, kdy sice fungují v principu stejně - uživatel může získat zvláštní práva - ale`n`tv případě licencí pouze na omezený čas, kdy pak buď musí licenci za poplatek obnovit, nebo si ji pořídit za plnou cenu znovu.`n</p>`n<p>`n`tHlavním smyslem licencí je umožnit na určitý vymezený čas uživateli přidělit nějaká zajímavá oprávnění.`n</p>`n<p>`n`tJelikož certifikát je vnímán jako exkluzivní zboží, některé základní vlastnosti aplikace jsou pro všechny uživatele řešené licencemi,`n`tjedná se tak například o přístup na Tržiště nebo do Laboratoře. Tyto jsou na rozdíl od standardní licence také neomezené.`n</p>"
$ws.Cells.Item(1305, 3).WrapText = $true
$ws.Cells.Item(1305, 3).Font.Size = 10

$ws.Rows.Item(1305).RowHeight = 141

$ws.Range("B1296").Select()
$excel.ActiveWindow.ScrollRow = 1279
$excel.ActiveWindow.ScrollColumn = 1
